$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows 11-13 (MuSCs sender rows shift out; table now has 9 data rows instead of 12)
$ws.Range("A11:T13").Delete()

# Overwrite rows 2-10 with the recalculated TPM-derived values
# Row 2: ECs -> FAPs
$ws.Cells.Item(2,1).Value2 = "ECs"
$ws.Cells.Item(2,2).Value2 = "Vegfc"
$ws.Cells.Item(2,3).Value2 = "Vipr2"
$ws.Cells.Item(2,4).Value2 = "FAPs"
$ws.Cells.Item(2,5).Value2 = 3
$ws.Cells.Item(2,6).Value2 = 1
$ws.Cells.Item(2,7).Value2 = 6.775549333333333
$ws.Cells.Item(2,8).Value2 = 20.326648
$ws.Cells.Item(2,9).Value2 = 0.5307754563424079
$ws.Cells.Item(2,10).Value2 = 0.5307754563424079
$ws.Cells.Item(2,11).Value2 = 3
$ws.Cells.Item(2,12).Value2 = 1
$ws.Cells.Item(2,13).Value2 = 6.804012333333334
$ws.Cells.Item(2,14).Value2 = 20.412037
$ws.Cells.Item(2,15).Value2 = 0.5269116569106099
$ws.Cells.Item(2,16).Value2 = 0.5269116569106099
$ws.Cells.Item(2,17).Value2 = 46.10092122910844
$ws.Cells.Item(2,18).Value2 = 414.908291061976
$ws.Cells.Item(2,19).Value2 = 0.2796717751488632
$ws.Cells.Item(2,20).Value2 = 0.2796717751488632

# Row 3: ECs -> MuSCs
$ws.Cells.Item(3,1).Value2 = "ECs"
$ws.Cells.Item(3,2).Value2 = "Vegfc"
$ws.Cells.Item(3,3).Value2 = "Vipr2"
$ws.Cells.Item(3,4).Value2 = "MuSCs"
$ws.Cells.Item(3,5).Value2 = 3
$ws.Cells.Item(3,6).Value2 = 1
$ws.Cells.Item(3,7).Value2 = 6.775549333333333
$ws.Cells.Item(3,8).Value2 = 20.326648
$ws.Cells.Item(3,9).Value2 = 0.5307754563424079
$ws.Cells.Item(3,10).Value2 = 0.5307754563424079
$ws.Cells.Item(3,11).Value2 = 3
$ws.Cells.Item(3,12).Value2 = 1
$ws.Cells.Item(3,13).Value2 = 6.095937333333333
$ws.Cells.Item(3,14).Value2 = 18.287812
$ws.Cells.Item(3,15).Value2 = 0.4720773983600821
$ws.Cells.Item(3,16).Value2 = 0.472077398360082
$ws.Cells.Item(3,17).Value2 = 41.30332413490844
$ws.Cells.Item(3,18).Value2 = 371.7299172141759
$ws.Cells.Item(3,19).Value2 = 0.2505670965435092
$ws.Cells.Item(3,20).Value2 = 0.2505670965435092

# Row 4: ECs -> Resolving-Mac
$ws.Cells.Item(4,1).Value2 = "ECs"
$ws.Cells.Item(4,2).Value2 = "Vegfc"
$ws.Cells.Item(4,3).Value2 = "Vipr2"
$ws.Cells.Item(4,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(4,5).Value2 = 3
$ws.Cells.Item(4,6).Value2 = 1
$ws.Cells.Item(4,7).Value2 = 6.775549333333333
$ws.Cells.Item(4,8).Value2 = 20.326648
$ws.Cells.Item(4,9).Value2 = 0.5307754563424079
$ws.Cells.Item(4,10).Value2 = 0.5307754563424079
$ws.Cells.Item(4,11).Value2 = 1
$ws.Cells.Item(4,12).Value2 = 0.3333333333333333
$ws.Cells.Item(4,13).Value2 = 0.01305433333333333
$ws.Cells.Item(4,14).Value2 = 0.039163
$ws.Cells.Item(4,15).Value2 = 0.001010944729308016
$ws.Cells.Item(4,16).Value2 = 0.001010944729308016
$ws.Cells.Item(4,17).Value2 = 0.08845027951377778
$ws.Cells.Item(4,18).Value2 = 0.796052515624
$ws.Cells.Item(4,19).Value2 = 0.0005365846500354144
$ws.Cells.Item(4,20).Value2 = 0.0005365846500354144

# Row 5: FAPs -> FAPs
$ws.Cells.Item(5,1).Value2 = "FAPs"
$ws.Cells.Item(5,2).Value2 = "Vegfc"
$ws.Cells.Item(5,3).Value2 = "Vipr2"
$ws.Cells.Item(5,4).Value2 = "FAPs"
$ws.Cells.Item(5,5).Value2 = 3
$ws.Cells.Item(5,6).Value2 = 1
$ws.Cells.Item(5,7).Value2 = 4.367310666666667
$ws.Cells.Item(5,8).Value2 = 13.101932
$ws.Cells.Item(5,9).Value2 = 0.3421215311185197
$ws.Cells.Item(5,10).Value2 = 0.3421215311185197
$ws.Cells.Item(5,11).Value2 = 3
$ws.Cells.Item(5,12).Value2 = 1
$ws.Cells.Item(5,13).Value2 = 6.804012333333334
$ws.Cells.Item(5,14).Value2 = 20.412037
$ws.Cells.Item(5,15).Value2 = 0.5269116569106099
$ws.Cells.Item(5,16).Value2 = 0.5269116569106099
$ws.Cells.Item(5,17).Value2 = 29.71523563949822
$ws.Cells.Item(5,18).Value2 = 267.437120755484
$ws.Cells.Item(5,19).Value2 = 0.180267822826454
$ws.Cells.Item(5,20).Value2 = 0.180267822826454

# Row 6: FAPs -> MuSCs
$ws.Cells.Item(6,1).Value2 = "FAPs"
$ws.Cells.Item(6,2).Value2 = "Vegfc"
$ws.Cells.Item(6,3).Value2 = "Vipr2"
$ws.Cells.Item(6,4).Value2 = "MuSCs"
$ws.Cells.Item(6,5).Value2 = 3
$ws.Cells.Item(6,6).Value2 = 1
$ws.Cells.Item(6,7).Value2 = 4.367310666666667
$ws.Cells.Item(6,8).Value2 = 13.101932
$ws.Cells.Item(6,9).Value2 = 0.3421215311185197
$ws.Cells.Item(6,10).Value2 = 0.3421215311185197
$ws.Cells.Item(6,11).Value2 = 3
$ws.Cells.Item(6,12).Value2 = 1
$ws.Cells.Item(6,13).Value2 = 6.095937333333333
$ws.Cells.Item(6,14).Value2 = 18.287812
$ws.Cells.Item(6,15).Value2 = 0.4720773983600821
$ws.Cells.Item(6,16).Value2 = 0.472077398360082
$ws.Cells.Item(6,17).Value2 = 26.62285213919822
$ws.Cells.Item(6,18).Value2 = 239.605669252784
$ws.Cells.Item(6,19).Value2 = 0.1615078423333987
$ws.Cells.Item(6,20).Value2 = 0.1615078423333986

# Row 7: FAPs -> Resolving-Mac
$ws.Cells.Item(7,1).Value2 = "FAPs"
$ws.Cells.Item(7,2).Value2 = "Vegfc"
$ws.Cells.Item(7,3).Value2 = "Vipr2"
$ws.Cells.Item(7,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(7,5).Value2 = 3
$ws.Cells.Item(7,6).Value2 = 1
$ws.Cells.Item(7,7).Value2 = 4.367310666666667
$ws.Cells.Item(7,8).Value2 = 13.101932
$ws.Cells.Item(7,9).Value2 = 0.3421215311185197
$ws.Cells.Item(7,10).Value2 = 0.3421215311185197
$ws.Cells.Item(7,11).Value2 = 1
$ws.Cells.Item(7,12).Value2 = 0.3333333333333333
$ws.Cells.Item(7,13).Value2 = 0.01305433333333333
$ws.Cells.Item(7,14).Value2 = 0.039163
$ws.Cells.Item(7,15).Value2 = 0.001010944729308016
$ws.Cells.Item(7,16).Value2 = 0.001010944729308016
$ws.Cells.Item(7,17).Value2 = 0.05701232921288889
$ws.Cells.Item(7,18).Value2 = 0.513110962916
$ws.Cells.Item(7,19).Value2 = 0.0003458659586670561
$ws.Cells.Item(7,20).Value2 = 0.0003458659586670561

# Row 8: MuSCs -> FAPs
$ws.Cells.Item(8,1).Value2 = "MuSCs"
$ws.Cells.Item(8,2).Value2 = "Vegfc"
$ws.Cells.Item(8,3).Value2 = "Vipr2"
$ws.Cells.Item(8,4).Value2 = "FAPs"
$ws.Cells.Item(8,5).Value2 = 3
$ws.Cells.Item(8,6).Value2 = 1
$ws.Cells.Item(8,7).Value2 = 1.622518
$ws.Cells.Item(8,8).Value2 = 4.867554
$ws.Cells.Item(8,9).Value2 = 0.1271030125390725
$ws.Cells.Item(8,10).Value2 = 0.1271030125390725
$ws.Cells.Item(8,11).Value2 = 3
$ws.Cells.Item(8,12).Value2 = 1
$ws.Cells.Item(8,13).Value2 = 6.804012333333334
$ws.Cells.Item(8,14).Value2 = 20.412037
$ws.Cells.Item(8,15).Value2 = 0.5269116569106099
$ws.Cells.Item(8,16).Value2 = 0.5269116569106099
$ws.Cells.Item(8,17).Value2 = 11.03963248305534
$ws.Cells.Item(8,18).Value2 = 99.35669234749801
$ws.Cells.Item(8,19).Value2 = 0.06697205893529272
$ws.Cells.Item(8,20).Value2 = 0.06697205893529272

# Row 9: MuSCs -> MuSCs
$ws.Cells.Item(9,1).Value2 = "MuSCs"
$ws.Cells.Item(9,2).Value2 = "Vegfc"
$ws.Cells.Item(9,3).Value2 = "Vipr2"
$ws.Cells.Item(9,4).Value2 = "MuSCs"
$ws.Cells.Item(9,5).Value2 = 3
$ws.Cells.Item(9,6).Value2 = 1
$ws.Cells.Item(9,7).Value2 = 1.622518
$ws.Cells.Item(9,8).Value2 = 4.867554
$ws.Cells.Item(9,9).Value2 = 0.1271030125390725
$ws.Cells.Item(9,10).Value2 = 0.1271030125390725
$ws.Cells.Item(9,11).Value2 = 3
$ws.Cells.Item(9,12).Value2 = 1
$ws.Cells.Item(9,13).Value2 = 6.095937333333333
$ws.Cells.Item(9,14).Value2 = 18.287812
$ws.Cells.Item(9,15).Value2 = 0.4720773983600821
$ws.Cells.Item(9,16).Value2 = 0.472077398360082
$ws.Cells.Item(9,17).Value2 = 9.890768050205333
$ws.Cells.Item(9,18).Value2 = 89.016912451848
$ws.Cells.Item(9,19).Value2 = 0.06000245948317423
$ws.Cells.Item(9,20).Value2 = 0.06000245948317422

# Row 10: MuSCs -> Resolving-Mac
$ws.Cells.Item(10,1).Value2 = "MuSCs"
$ws.Cells.Item(10,2).Value2 = "Vegfc"
$ws.Cells.Item(10,3).Value2 = "Vipr2"
$ws.Cells.Item(10,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(10,5).Value2 = 3
$ws.Cells.Item(10,6).Value2 = 1
$ws.Cells.Item(10,7).Value2 = 1.622518
$ws.Cells.Item(10,8).Value2 = 4.867554
$ws.Cells.Item(10,9).Value2 = 0.1271030125390725
$ws.Cells.Item(10,10).Value2 = 0.1271030125390725
$ws.Cells.Item(10,11).Value2 = 1
$ws.Cells.Item(10,12).Value2 = 0.3333333333333333
$ws.Cells.Item(10,13).Value2 = 0.01305433333333333
$ws.Cells.Item(10,14).Value2 = 0.039163
$ws.Cells.Item(10,15).Value2 = 0.001010944729308016
$ws.Cells.Item(10,16).Value2 = 0.001010944729308016
$ws.Cells.Item(10,17).Value2 = 0.02118089081133334
$ws.Cells.Item(10,18).Value2 = 0.190628017302
$ws.Cells.Item(10,19).Value2 = 0.0001284941206055461
$ws.Cells.Item(10,20).Value2 = 0.0001284941206055461
